$p = $ppt.ActivePresentation

# --- 1. Slide 5 table: switch to a different built-in table style -----------
$tableShape = $null
$s5 = $p.Slides.Item(5)
for ($i = 1; $i -le $s5.Shapes.Count; $i++) {
    $sh = $s5.Shapes.Item($i)
    if ($sh.HasTable) {
        $tableShape = $sh
        break
    }
}
if ($tableShape -ne $null) {
    $tableShape.Table.ApplyStyle("{1D11E14F-EDD9-4D33-A21A-1FB770C0110F}", $true)
}

# --- 2. Theme colours: swap in the "Office" palette --------------------------
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order)
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$theme = $p.SlideMaster.Theme
for ($i = 1; $i -le 12; $i++) {
    $hex = $officeColors[$i - 1]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    $comRgb = $r -bor ($g -shl 8) -bor ($b -shl 16)
    $theme.ThemeColorScheme.Colors($i).RGB = $comRgb
}
